$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B80").Value = 7
$ws.Range("C80").Value = "bugs fixed, other bugs added, tiles redraw"

$ws.Range("B81").Value = 6
$ws.Range("C81").Value = "level design"

$ws.Range("C80").Select()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
